# de_dg_gw_upper.xlsx - "add rural communities interactions back in, fix de_dg files (not run yet)"
#
# The "Legislature" row (row 12) is removed from the interactions table;
# everything below it shifts up by one row, and the now-unused
# "Legislature" shared string is dropped automatically by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire row for "Legislature" (row 12). Excel will shift all
# subsequent rows up by one and clean up the now-unused shared string.
[void]($ws.Rows(12).Delete())

# Restore the view/selection state recorded for this edit: the row that
# used to hold "Legislature" (now "Friant Water Authority") is selected...
[void]($ws.Rows(12).Select())

# ...and the window is scrolled so row 9 is at the top.
[void]($excel.ActiveWindow.ScrollRow = 9)
[void]($excel.ActiveWindow.ScrollColumn = 1)
